# Text-only edit for slide 1 (images are left untouched - harder to automate).
#
#  - Title placeholder:         "Test" -> "Cjfshbvfhjsbfhsb"
#  - Body placeholder (idx=18): "Content would be here Why bullet points"
#                                -> "kdsnxzbfjhsdbfdsbfu"
#  - Body placeholder (idx=19): "- God" / "- John" / "- Jesus" bullets turned
#                                into a long run of single-letter bullets.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Title ("Test" -> "Cjfshbvfhjsbfhsb") ------------------------------
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleLen = $titleRange.Text.Length
$titleRange.Characters(1, $titleLen).Text = "Cjfshbvfhjsbfhsb"

# ---- Body placeholder idx=18 (content placeholder) ---------------------
$contentShape = $s.Shapes.Item(2)
$contentRange = $contentShape.TextFrame.TextRange
$contentRange.Text = "kdsnxzbfjhsdbfdsbfu"
$contentRange.LanguageID = "en-GB"

# ---- Body placeholder idx=19 ("People involved:" bullets) --------------
$peopleShape = $s.Shapes.Item(3)
$peopleRange = $peopleShape.TextFrame.TextRange

# Replace the existing 3 bullet paragraphs ("- God"/"- John"/"- Jesus")
# in place so they keep their paragraph/run formatting (lvl="1" etc.).
$peopleRange.Characters(18, 5).Text = "- C"
$peopleRange.Characters(22, 6).Text = "- j"
$peopleRange.Characters(26, 7).Text = "- f"

# Append the remaining bullets as brand new paragraphs at the same
# (inherited) indent level.
$peopleRange.InsertAfter("`r- s`r- h`r- b`r- v`r- f`r- h`r- j`r- s`r- b`r- f`r- h`r- s`r- b") | Out-Null
